$d = $word.ActiveDocument

# Paragraph 1: "Date: 23 September 2023" -> "Date: 25 September 2023"
$r1 = $d.Paragraphs(1).Range
$r1.Find.Execute("Date: 23 September 2023", $false, $false, $false, $false, $false, $true, 1, $false, "Date: 25 September 2023", 2)

# Paragraph 2: "Start Time: 12:46 PM" -> "Start Time: 6:00 PM"
$r2 = $d.Paragraphs(2).Range
$r2.Find.Execute("Start Time: 12:46 PM", $false, $false, $false, $false, $false, $true, 1, $false, "Start Time: 6:00 PM", 2)

# Paragraph 3: "End Time:  6:46 PM" -> "End Time:  8:30 PM"
$r3 = $d.Paragraphs(3).Range
$r3.Find.Execute("End Time:  6:46 PM", $false, $false, $false, $false, $false, $true, 1, $false, "End Time:  8:30 PM", 2)

# Paragraph 4: "Total Time: 19 hour 23 Minutes" -> "Total Time: 23 hour 23 Minutes"
$r4 = $d.Paragraphs(4).Range
$r4.Find.Execute("Total Time: 19 hour 23 Minutes", $false, $false, $false, $false, $false, $true, 1, $false, "Total Time: 23 hour 23 Minutes", 2)
